# Auto commit at 2025-10-03  9:04:28.65
# Append the two new daily rows (四方坪站 / 高岭站) for 2025-10-02
# and move the active selection to H8, matching the source workbook update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: 四方坪站 for date 45932 (2025-10-02)
$ws.Range("A4").Value = 45932
$ws.Range("B4").Value = "四方坪站"
$ws.Range("C4").Value = 9431.59
$ws.Range("D4").Value = 7990.26
$ws.Range("E4").Value = 3228.21
$ws.Range("F4").Value = 383

# Row 5: 高岭站 for date 45932 (2025-10-02)
$ws.Range("A5").Value = 45932
$ws.Range("B5").Value = "高岭站"
$ws.Range("C5").Value = 4692.01
$ws.Range("D5").Value = 3750.2
$ws.Range("E5").Value = 1213.94
$ws.Range("F5").Value = 153

# Update the saved selection/active cell, as seen in the diff.
$ws.Range("H8").Select()
